$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    # Plain whole-document replace. Used where the paragraph either has no
    # sibling empty run to worry about, or its text run already carries
    # explicit run formatting (which keeps the engine from collapsing the
    # paragraph's run list during the rewrite).
    $d.Content.Find.Execute(
        $findText, $true, $false, $false, $false, $false,
        $true, 1, $false, $replaceText, 2) | Out-Null
}

function Replace-Text-KeepEmptyRun($findText, $replaceText) {
    # Some paragraphs in this document are `<w:r/><w:r><w:t>...</w:t></w:r>`
    # - an empty leading run followed by the text run. When the text run's
    # Find/Replace rewrites that paragraph, the engine drops the empty
    # sibling run unless the edited run already carries run-level formatting
    # (<w:rPr>). Briefly toggling Bold on/off around the text replacement
    # gives the run an <w:rPr> for the rewrite and keeps the empty run
    # intact, then clears the formatting back off again so the run ends up
    # plain, matching the target text exactly.
    $r = $d.Content
    $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $r.Bold = 1
    $d.Content.Find.Execute(
        $findText, $true, $false, $false, $false, $false,
        $true, 1, $false, $replaceText, 2) | Out-Null
    $r2 = $d.Content
    $r2.Find.Execute($replaceText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $r2.Bold = 0
}

# Title / heading text changed throughout the document (both occurrences are
# identical before the edit and identical after, so a document-wide replace
# is safe here).
Replace-Text "Play Mermaid Reef Free Today - Slot Game Review" "Play Mermaid Reef for Free"

# "What we like" bullet list
Replace-Text-KeepEmptyRun "Flexible betting range" "Flexible betting options"
Replace-Text-KeepEmptyRun "Autoplay feature with no limits" "Frequent winnings"
Replace-Text-KeepEmptyRun "Bonus game with multipliers up to 20x" "Beautiful visual design"

# "What we don't like" bullet list. The first bullet's new text is the
# second bullet's current (pre-edit) text, and the second bullet gets a new
# value entirely, so handle the second (currently unique) occurrence first.
#
# Reuse a single Range ($r) for both Find.Execute calls: after a successful
# match, Word's Find collapses the range to the match and the next
# Find.Execute on that same Range continues searching forward from there.
# That lets us walk past the first (still unmodified) bullet and land on the
# second occurrence of "Lack of big special effects" without manually
# recomputing Start/End offsets (manual Start/End assignment on a Range has
# its own side effects in this runtime, so it's avoided).
$r = $d.Content
$r.Find.Execute("Limited number of paylines", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Bold = 1
$r.Find.Execute(
    "Lack of big special effects", $true, $false, $false, $false, $false,
    $true, 1, $false, "Limited bonus game activation", 2) | Out-Null
$r2 = $d.Content
$r2.Find.Execute("Limited bonus game activation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Bold = 0

# Now replace the first bullet using a fresh full-document range.
Replace-Text-KeepEmptyRun "Limited number of paylines" "Lack of big special effects"

# Meta title / description runs near the end of the document (already bold /
# italic, so their <w:rPr> survives the rewrite naturally).
Replace-Text "Discover the beautiful Abaco Islands with the Mermaid Reef slot game from ReelPlay. Play and enjoy generous bonus games and high volatility for big payouts." "Read our review of Mermaid Reef and play this slot game for free."
